# Check supp one marks
#
# The four data-bound cells (StudentNo/StudName/Exam_Marks/CATs_Marks) and the
# repeating-section container sdt get a Title (w:alias) + Tag (w:tag) set, as
# happens when the content control's "Properties" are filled in via the
# Developer tab in Word. Saving also normalizes the stray trailing space
# Word/NAV had left inside the dataBinding prefixMappings attribute.
#
# Because this document's bound cells live inside a repeating-section
# template row, the simplified Word object model here does not walk them as
# ordinary ContentControls/Paragraphs, so we round-trip the whole document
# through Range.WordOpenXML / Range.InsertXML to touch the underlying SDTs.

$d = $word.ActiveDocument
$full = $d.Range(0, $d.Characters.Count)
$xml = $full.WordOpenXML

# 1) Drop the trailing space Word NAV left inside every dataBinding's
#    prefixMappings="...' " attribute value.
$xml = $xml.Replace("Students_Marks_Upload/50820/' `"", "Students_Marks_Upload/50820/'`"")

# 2) Repeating-section container sdt: add alias/tag after w15:repeatingSection.
$xml = $xml.Replace(
  "<w15:repeatingSection/></w:sdtPr>",
  "<w15:repeatingSection/><w:alias w:val=`"#Nav: /ACAStudentUnits`"/><w:tag w:val=`"#Nav: Students_Marks_Upload/50820`"/></w:sdtPr>"
)

# 3) The four bound-cell sdt's: add alias/tag after <w:text/>, each keyed to
#    the element name referenced by its own xpath so every cell gets the
#    right alias.
$xml = $xml.Replace(
  "ACAStudentUnits[1]/ns0:StudentNo[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/></w:sdtPr>",
  "ACAStudentUnits[1]/ns0:StudentNo[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/><w:alias w:val=`"#Nav: /ACAStudentUnits/StudentNo`"/><w:tag w:val=`"#Nav: Students_Marks_Upload/50820`"/></w:sdtPr>"
)
$xml = $xml.Replace(
  "ACAStudentUnits[1]/ns0:StudName[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/></w:sdtPr>",
  "ACAStudentUnits[1]/ns0:StudName[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/><w:alias w:val=`"#Nav: /ACAStudentUnits/StudName`"/><w:tag w:val=`"#Nav: Students_Marks_Upload/50820`"/></w:sdtPr>"
)
$xml = $xml.Replace(
  "ACAStudentUnits[1]/ns0:Exam_Marks[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/></w:sdtPr>",
  "ACAStudentUnits[1]/ns0:Exam_Marks[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/><w:alias w:val=`"#Nav: /ACAStudentUnits/Exam_Marks`"/><w:tag w:val=`"#Nav: Students_Marks_Upload/50820`"/></w:sdtPr>"
)
$xml = $xml.Replace(
  "ACAStudentUnits[1]/ns0:CATs_Marks[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/></w:sdtPr>",
  "ACAStudentUnits[1]/ns0:CATs_Marks[1]`" w:storeItemID=`"{7F66C27A-F9FF-4E01-851E-010AB93FF16B}`"/><w:text/><w:alias w:val=`"#Nav: /ACAStudentUnits/CATs_Marks`"/><w:tag w:val=`"#Nav: Students_Marks_Upload/50820`"/></w:sdtPr>"
)

$full.InsertXML($xml)
